# "Checked grapfics for bubblesort"
# Update the BubbleSort results/row data, clear the other algorithms'
# "Random Numbers" rows, fill in the two missing BubbleSort timings for
# the "Numbers in Order" / "Numbers in Reverse Order" tables, reposition
# the chart, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Random Numbers" table (rows 2-9) ---------------------------------
# BubbleSort (row 3) timings were re-measured.
$ws.Range("D3").Value = 7544
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 5

# The other algorithms (HeapSort, InsertionSort, MergeSort, QuickSort -
# rows 4-7) haven't been benchmarked yet, so their results are cleared.
$ws.Range("B4:G7").ClearContents()

# --- "Numbers in Order" table (rows 11-18) ------------------------------
# BubbleSort (row 12) timing for the 50000-element case.
$ws.Range("D12").Value = 3436

# --- "Numbers in Reverse Order" table (rows 20-27) ----------------------
# BubbleSort (row 21) timing for the 50000-element case.
$ws.Range("D21").Value = 6108

# --- Reposition the chart -----------------------------------------------
$co = $ws.ChartObjects(1)
$co.Left = 543.56640625
$co.Top = 40.5
$co.Width = 638.3125
$co.Height = 418.4999212598425

# --- Active selection ----------------------------------------------------
$ws.Range("G4").Select()
